$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.360.92'
$ws.Range("E2").Value = '  +2.81%  '

$ws.Range("D3").Value = '2.064.15'
$ws.Range("E3").Value = '  +4.48%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.75'
$ws.Range("E5").Value = '  +1.73%  '

$ws.Range("E6").Value = '  +3.49%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.05'
$ws.Range("E7").Value = '  +6.69%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.383'
$ws.Range("E9").Value = '  +3.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.94'
$ws.Range("E10").Value = '  -1.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0762'
$ws.Range("E11").Value = '  +2.12%  '

$ws.Range("E12").Value = '  +3.62%  '

$ws.Range("D13").Value = '2.367.49'
$ws.Range("E13").Value = '  +4.42%  '

$ws.Range("E14").Value = '  +4.11%  '

$ws.Range("E15").Value = '  +5.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.779'
$ws.Range("E16").Value = '  +4.12%  '

$ws.Range("E17").Value = '  +3.10%  '

$ws.Range("D18").Value = '2.065.18'
$ws.Range("E18").Value = '  +4.59%  '

$ws.Range("D19").Value = '37.441.39'
$ws.Range("E19").Value = '  +3.15%  '

$ws.Range("E20").Value = '  +17.67%  '

$ws.Range("E21").Value = '  +2.49%  '

$ws.Range("D22").Value = '0.0₃0816'
$ws.Range("E22").Value = '  +1.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '228.18'
$ws.Range("E23").Value = '  +3.22%  '

$ws.Range("E24").Value = '  +0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("E25").Value = '  +3.36%  '

$ws.Range("E26").Value = '  +1.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.88'
$ws.Range("E27").Value = '  +2.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.48'
$ws.Range("E28").Value = '  +11.96%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.87'
$ws.Range("E29").Value = '  +3.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.17'
$ws.Range("E30").Value = '  +2.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.127'
$ws.Range("E31").Value = '  +2.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.118'
$ws.Range("E32").Value = '  +2.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.52'
$ws.Range("E33").Value = '  +4.09%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.59'
$ws.Range("E34").Value = '  +12.72%  '

$ws.Range("E35").Value = '  +2.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.49'
$ws.Range("E36").Value = '  +6.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.42'

$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.92'
$ws.Range("E39").Value = '  +11.68%  '

$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.78'
$ws.Range("E40").Value = '  +0.52%  '

$ws.Range("B41").Value = 'FTXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.60'
$ws.Range("E41").Value = '  +30.50%  '

$ws.Range("B42").Value = 'Cronos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0991'
$ws.Range("E42").Value = '  +10.63%  '

$ws.Range("E43").Value = '  -1.87%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.33'
$ws.Range("E44").Value = '  +10.74%  '

$ws.Range("D45").Value = '1.477.29'
$ws.Range("E45").Value = '  +1.76%  '

$ws.Range("E46").Value = '  +8.04%  '

$ws.Range("E47").Value = '  +5.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.90'
$ws.Range("E48").Value = '  +7.48%  '

$ws.Range("E49").Value = '  +3.91%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.21'
$ws.Range("E50").Value = '  +6.87%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.93'
$ws.Range("E51").Value = '  +2.04%  '
